# Minor fix in TSP.
# Update the "Fitness" column (C) values for generations 0-10 (rows 2-12)
# on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = 3870
$ws.Range("C3").Value  = 3870
$ws.Range("C4").Value  = 3911
$ws.Range("C5").Value  = 4014
$ws.Range("C6").Value  = 4035
$ws.Range("C7").Value  = 4035
$ws.Range("C8").Value  = 4035
$ws.Range("C9").Value  = 4035
$ws.Range("C10").Value = 4051
$ws.Range("C11").Value = 4669
$ws.Range("C12").Value = 4669
